# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "002" -> "001" (keep as text so the leading zero survives)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").Style = "Normal"

# Dates (stored as text strings)
$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Numeric figures
$ws.Range("O2").Value = 16759911.2
$ws.Range("P2").Value = 422647488.68
$ws.Range("Q2").Value = 392756729.05
$ws.Range("R2").Value = 40.3130994325
$ws.Range("S2").Value = 305106899.47
$ws.Range("T2").Value = 305106899.47
$ws.Range("U2").Value = 39.1694943121
$ws.Range("V2").Value = 16243638.65
$ws.Range("W2").Value = 50217349.71
$ws.Range("X2").Value = 3967598.48
$ws.Range("Y2").Value = 19444531.92
$ws.Range("Z2").Value = 20846338.67
$ws.Range("AA2").Value = 4203424.3

$ws.Range("AG2").Value = 3094261.29

$ws.Range("AP2").Value = 34.4340118284
$ws.Range("AQ2").Value = -43.599449058499
$ws.Range("AR2").Value = -46.394312217299
$ws.Range("AS2").Value = 15321709.33
$ws.Range("AT2").Value = -49.213923251622
